$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 77 ("Macroferia Regional de
# Talca" / Mango data block). This shifts the existing rows 77-142 down to
# 78-143 (matching the dimension change from A1:T142 to A1:T143) while
# leaving row 76 and everything above untouched.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new data record.
$ws.Range("A77").Value = 5
$ws.Range("B77").Value = "Macroferia Regional de Talca"
$ws.Range("C77").Value = "Maule"
$ws.Range("D77").Value = 44790
$ws.Range("E77").Value = 7
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100108
$ws.Range("H77").Value = "Tropicales y subtropicales"
$ws.Range("I77").Value = 100108002
$ws.Range("J77").Value = "Mango"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 228
$ws.Range("N77").Value = 10000
$ws.Range("O77").Value = 10000
$ws.Range("P77").Value = 10000
$ws.Range("Q77").Value = "$/bandeja 4 kilos"
$ws.Range("R77").Value = "Brasil"
$ws.Range("S77").Value = 2500
$ws.Range("T77").Value = 4
